$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.861.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.188.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.21%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.531"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.94%  "
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.430"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.737.54"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.890.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.186.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "366.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.520"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.32%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0876"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.03%  "
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.79%  "
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "155.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.793.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.57%  "
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.91"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0699"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("E42").Value = "  +4.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.715"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("E45").Value = "  +2.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.227.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.986"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.797"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.04%  "
